# "last of the sliding window neetcodes"
#
# Fill in the missing STATUS / NOTES for "Minimum Window Substring" (row 20)
# and append a brand-new tracker row for "Sliding Window Maximum" (row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leetcode")

$yellow = 65535  # RGB(255,255,0) -> used for the "STRUGGLED" highlight rows

# --- Row 20: "Minimum Window Substring" was missing STATUS + NOTES ---
$ws.Range("G20").Value = "STRUGGLED"
$ws.Range("I20").Value = "Another tough one."

# This row now joins the "struggled" rows, which get a yellow row highlight,
# a left-aligned IDENTIFIER cell and a wrap-texted NOTES cell.
$ws.Rows.Item(20).Interior.Color = $yellow
$ws.Range("I20").WrapText = $true

# --- Row 21: brand-new entry for "Sliding Window Maximum" ---
$ws.Range("A21").Value = "Leetcode"
$ws.Range("B21").Value = 239
$ws.Range("C21").Value = "Sliding Window Maximum"
$ws.Range("D21").Value = "Arrays, Sliding Window, Queue"
$ws.Range("E21").Value = "Hard"
$ws.Range("F21").Value = "Neetcode 150"
$ws.Range("G21").Value = "STRUGGLED"
# Force the date to be stored as literal text (matching every other LAST
# SOLVED cell), not auto-converted to a date serial number.
$ws.Range("H21").Value = "'07/06/2025"
$ws.Range("H21").ClearFormats()
$ws.Range("I21").Value = "Using a special kind of Queue called a Deque in Python. We're using a Monotonically Decreasing Queue."

$ws.Rows.Item(21).Interior.Color = $yellow
$ws.Range("I21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 30

# The "Top K Frequent Elements" row (6) also grows to a two-line row height
# once the notes column wraps at the new column width.
$ws.Rows.Item(6).RowHeight = 30

# --- View / selection bookkeeping ---
$ws.Range("F21").Select()

$wb.Save()
